$d = $word.ActiveDocument

# 1. Remove the paragraph-level rFonts/pPr formatting from the paragraph
#    "玩家站在梯子附近时（很近）按W角色应会自动走到梯子附近爬上去..."
#    by clearing its paragraph format run properties (rPr inside pPr).
foreach ($para in $d.Paragraphs) {
    $text = $para.Range.Text
    if ($text -match "玩家站在梯子附近时") {
        $para.Range.ParagraphFormat.Reset()
    }
}

# 2. Add a new run with text "在梯子上不能左右移动" into the last paragraph
#    (the one that only contains the _GoBack bookmark), right before the
#    bookmark.
$lastPara = $d.Paragraphs.Last
$insertRange = $lastPara.Range
$insertRange.Collapse(1)
$insertRange.Font.Reset()
$insertRange.InsertBefore("在梯子上不能左右移动")
